$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Categories").Name = "Hierarchies"
$wb.Worksheets.Item("CategoriesMapping").Name = "HierarchiesMapping"

$ws1 = $wb.Worksheets.Item("Hierarchies")
$ws1.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws1.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

$ws2 = $wb.Worksheets.Item("HierarchiesMapping")
$ws2.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws2.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
